# Applies two changes described by the commit diff:
#  1. Delete the leading empty paragraph at the very start of the document
#     (an empty paragraph that precedes the "<page>074r</page>" paragraph).
#  2. Fix a spelling typo later in the document: "asses" -> "assés"
#     (only the lone "e" run changes to "é", leaving the surrounding runs
#     and their formatting untouched).

$d = $word.ActiveDocument

# --- Change 1: remove the empty first paragraph -----------------------
# Paragraph.Range.Text includes the trailing paragraph mark, so an
# "empty" paragraph has a Range.Text of just that one mark character.
$firstPara = $d.Paragraphs(1)
if ($firstPara.Range.Text.Length -le 1) {
    $firstPara.Range.Delete()
}

# --- Change 2: "asses" -> "assés" (single-character fix) --------------
# Locate the word "asses" (only occurrence in the document, immediately
# followed by " cuyt") and narrow the range down to just the "e"
# character so only that run's text is touched (matches the diff, which
# only rewrites the <w:t> content of that one run).
$found = $d.Content
$ok = $found.Find.Execute("asses cuyt", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok) {
    $wordStart = $found.Start
    $charRange = $d.Range($wordStart + 3, $wordStart + 4)
    $charRange.Text = "é"
}
